# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1038.5
$ws.Range("I4").Value = 1046.2
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1046.2
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -932.2
$ws.Range("N4").Value = -1228

$ws.Range("H8").Value = 1032.8572
$ws.Range("I8").Value = 743.3333
$ws.Range("J8").Value = 1250
$ws.Range("K8").Value = 2229.9999
$ws.Range("L8").Value = 3750
$ws.Range("M8").Value = -2090.9999
$ws.Range("N8").Value = -4028

$ws.Range("H18").Value = 933.7
$ws.Range("I18").Value = 771.4666999999999
$ws.Range("J18").Value = 1420.4
$ws.Range("K18").Value = 771.4666999999999
$ws.Range("L18").Value = 1420.4
$ws.Range("M18").Value = -487.4666999999999
$ws.Range("N18").Value = -1988.4

$ws.Range("H32").Value = 2644.5
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 2887.5625
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 2887.5625
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -3539.5625

$ws.Range("H41").Value = 314.05554
$ws.Range("I41").Value = 593.6667
$ws.Range("J41").Value = 258.13333
$ws.Range("K41").Value = 593.6667
$ws.Range("L41").Value = 258.13333
$ws.Range("M41").Value = -153.6667
$ws.Range("N41").Value = -1138.13333

$ws.Range("H74").Value = 3933.158
$ws.Range("I74").Value = 3838.25
$ws.Range("J74").Value = 4002.182
$ws.Range("K74").Value = 3838.25
$ws.Range("L74").Value = 4002.182
$ws.Range("M74").Value = -2902.25
$ws.Range("N74").Value = -5874.182

$ws.Range("H77").Value = 3933.158
$ws.Range("I77").Value = 3838.25
$ws.Range("J77").Value = 4002.182
$ws.Range("K77").Value = 19191.25
$ws.Range("L77").Value = 20010.91
$ws.Range("M77").Value = -14511.25
$ws.Range("N77").Value = -29370.91

$ws.Range("H98").Value = 2884.75
$ws.Range("I98").Value = 3244.5
$ws.Range("J98").Value = 2525
$ws.Range("K98").Value = 3244.5
$ws.Range("L98").Value = 2525
$ws.Range("M98").Value = -1746.5
$ws.Range("N98").Value = -5521

$ws.Range("H122").Value = 2884.75
$ws.Range("I122").Value = 3244.5
$ws.Range("J122").Value = 2525
$ws.Range("K122").Value = 9733.5
$ws.Range("L122").Value = 7575
$ws.Range("M122").Value = -7283.5
$ws.Range("N122").Value = -12475

$ws.Range("H132").Value = 3504025
$ws.Range("I132").Value = 6333.875
$ws.Range("J132").Value = 8167613.5
$ws.Range("K132").Value = 19001.625
$ws.Range("L132").Value = 24502840.5
$ws.Range("M132").Value = -16471.625
$ws.Range("N132").Value = -24507900.5

$ws.Range("H137").Value = 2780224.8
$ws.Range("I137").Value = 4547460
$ws.Range("K137").Value = 13642380
$ws.Range("M137").Value = -13639830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3008.0527
$ws.Range("I63").Value = 2765.6667
$ws.Range("J63").Value = 3119.923
$ws.Range("K63").Value = 2765.6667
$ws.Range("L63").Value = 3119.923
$ws.Range("M63").Value = -2079.6667
$ws.Range("N63").Value = -4491.923

$ws.Range("H66").Value = 3008.0527
$ws.Range("I66").Value = 2765.6667
$ws.Range("J66").Value = 3119.923
$ws.Range("K66").Value = 13828.3335
$ws.Range("L66").Value = 15599.615
$ws.Range("M66").Value = -10396.3335
$ws.Range("N66").Value = -22463.615

$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800

$ws.Range("H138").Value = 44000
$ws.Range("J138").Value = 44000
$ws.Range("L138").Value = 44000
$ws.Range("N138").Value = -54280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1915.1305
$ws.Range("I31").Value = 1528.8422
$ws.Range("K31").Value = 1528.8422
$ws.Range("M31").Value = -1233.8422

$ws.Range("H34").Value = 1915.1305
$ws.Range("I34").Value = 1528.8422
$ws.Range("K34").Value = 1528.8422
$ws.Range("M34").Value = -1326.8422

$ws.Range("H105").Value = 2062.2
$ws.Range("I105").Value = 836.6667
$ws.Range("J105").Value = 3900.5
$ws.Range("K105").Value = 836.6667
$ws.Range("L105").Value = 3900.5
$ws.Range("M105").Value = 910.3333
$ws.Range("N105").Value = -7394.5

$ws.Range("H134").Value = 64616.41
$ws.Range("I134").Value = 908.0909
$ws.Range("K134").Value = 2724.2727
$ws.Range("M134").Value = -189.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 34138.3
$ws.Range("I5").Value = 53150.21
$ws.Range("K5").Value = 159450.63
$ws.Range("M5").Value = -159338.63

$ws.Range("H122").Value = 1291.9166
$ws.Range("J122").Value = 1544.3334
$ws.Range("L122").Value = 13899.0006
$ws.Range("N122").Value = -18799.0006

$ws.Range("H131").Value = 955.8108
$ws.Range("J131").Value = 979
$ws.Range("L131").Value = 2937
$ws.Range("N131").Value = -13017

$ws.Range("H132").Value = 1479.7222
$ws.Range("I132").Value = 1781
$ws.Range("K132").Value = 16029
$ws.Range("M132").Value = -13499

$ws.Range("H135").Value = 34138.3
$ws.Range("I135").Value = 53150.21
$ws.Range("K135").Value = 478351.89
$ws.Range("M135").Value = -475816.89

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2473.3333
$ws.Range("I122").Value = 2181.5454
$ws.Range("J122").Value = 2931.8572
$ws.Range("K122").Value = 6544.6362
$ws.Range("L122").Value = 8795.571599999999
$ws.Range("M122").Value = -4094.6362
$ws.Range("N122").Value = -13695.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1001.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1001.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1001.5
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -1377.5

$ws.Range("H68").Value = 1589.9524
$ws.Range("I68").Value = 1474.625
$ws.Range("J68").Value = 1959
$ws.Range("K68").Value = 1474.625
$ws.Range("L68").Value = 1959
$ws.Range("M68").Value = -725.625
$ws.Range("N68").Value = -3457

$ws.Range("H71").Value = 1589.9524
$ws.Range("I71").Value = 1474.625
$ws.Range("J71").Value = 1959
$ws.Range("K71").Value = 7373.125
$ws.Range("L71").Value = 9795
$ws.Range("M71").Value = -3629.125
$ws.Range("N71").Value = -17283

$ws.Range("H122").Value = 3855.889
$ws.Range("I122").Value = 3260.8
$ws.Range("J122").Value = 4599.75
$ws.Range("K122").Value = 9782.400000000001
$ws.Range("L122").Value = 13799.25
$ws.Range("M122").Value = -7332.400000000001
$ws.Range("N122").Value = -18699.25

$ws.Range("H136").Value = 290089.1
$ws.Range("I136").Value = 168833.5
$ws.Range("J136").Value = 471972.5
$ws.Range("K136").Value = 506500.5
$ws.Range("L136").Value = 1415917.5
$ws.Range("M136").Value = -503950.5
$ws.Range("N136").Value = -1421017.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2053.3333
$ws.Range("I81").Value = 1556.2
$ws.Range("J81").Value = 2674.75
$ws.Range("K81").Value = 3112.4
$ws.Range("L81").Value = 5349.5
$ws.Range("M81").Value = -2051.4
$ws.Range("N81").Value = -7471.5

$ws.Range("H84").Value = 2053.3333
$ws.Range("I84").Value = 1556.2
$ws.Range("J84").Value = 2674.75
$ws.Range("K84").Value = 15562
$ws.Range("L84").Value = 26747.5
$ws.Range("M84").Value = -10258
$ws.Range("N84").Value = -37355.5
